$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price (column D) values are plain numbers (e.g. "583.92") that Excel
# would otherwise auto-convert to a numeric type. Force those specific cells
# to text first so the textual formatting (e.g. trailing zeros) is preserved,
# exactly like the original inline-string cells.
$numericPriceCells = @("D5","D6","D8","D9","D11","D13","D15","D20","D23","D24","D28","D30","D32","D34","D37","D38","D39","D43","D44")
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = '68.270.61'
$ws.Cells.Item(2, 5).Value = '  +0.97%  '
$ws.Cells.Item(3, 4).Value = '3.349.65'
$ws.Cells.Item(3, 5).Value = '  +0.34%  '
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$ws.Cells.Item(5, 4).Value = '583.92'
$ws.Cells.Item(5, 5).Value = '  +0.32%  '
$ws.Cells.Item(6, 4).Value = '177.23'
$ws.Cells.Item(6, 5).Value = '  +0.62%  '
$ws.Cells.Item(8, 4).Value = '0.590'
$ws.Cells.Item(8, 5).Value = '  -0.15%  '
$ws.Cells.Item(9, 4).Value = '0.184'
$ws.Cells.Item(9, 5).Value = '  +2.48%  '
$ws.Cells.Item(10, 5).Value = '  +0.52%  '
$ws.Cells.Item(11, 4).Value = '48.12'
$ws.Cells.Item(11, 5).Value = '  +5.82%  '
$ws.Cells.Item(12, 5).Value = '  +1.42%  '
$ws.Cells.Item(13, 4).Value = '687.42'
$ws.Cells.Item(13, 5).Value = '  +3.00%  '
$ws.Cells.Item(14, 4).Value = '3.896.75'
$ws.Cells.Item(14, 5).Value = '  +0.45%  '
$ws.Cells.Item(15, 4).Value = '8.43'
$ws.Cells.Item(15, 5).Value = '  -0.06%  '
$ws.Cells.Item(16, 4).Value = '68.304.91'
$ws.Cells.Item(16, 5).Value = '  +0.75%  '
$ws.Cells.Item(17, 5).Value = '  +1.19%  '
$ws.Cells.Item(18, 4).Value = '3.347.45'
$ws.Cells.Item(18, 5).Value = '  +0.23%  '
$ws.Cells.Item(19, 5).Value = '  +0.16%  '
$ws.Cells.Item(20, 4).Value = '11.20'
$ws.Cells.Item(20, 5).Value = '  +2.07%  '
$ws.Cells.Item(21, 5).Value = '  +0.27%  '
$ws.Cells.Item(22, 5).Value = '  -0.51%  '
$ws.Cells.Item(23, 4).Value = '17.01'
$ws.Cells.Item(23, 5).Value = '  -0.61%  '
$ws.Cells.Item(24, 4).Value = '100.48'
$ws.Cells.Item(24, 5).Value = '  +0.98%  '
$ws.Cells.Item(25, 5).Value = '  +1.31%  '
$ws.Cells.Item(26, 5).Value = '  +0.95%  '
$ws.Cells.Item(27, 5).Value = '  +1.95%  '
$ws.Cells.Item(28, 4).Value = '32.98'
$ws.Cells.Item(28, 5).Value = '  -2.22%  '
$ws.Cells.Item(29, 5).Value = '  +0.79%  '
$ws.Cells.Item(30, 4).Value = '6.95'
$ws.Cells.Item(30, 5).Value = '  -6.49%  '
$ws.Cells.Item(31, 5).Value = '  +0.77%  '
$ws.Cells.Item(32, 4).Value = '555.97'
$ws.Cells.Item(32, 5).Value = '  -4.09%  '
$ws.Cells.Item(33, 5).Value = '  +0.82%  '
$ws.Cells.Item(34, 4).Value = '58.05'
$ws.Cells.Item(34, 5).Value = '  +2.46%  '
$ws.Cells.Item(35, 5).Value = '  +0.01%  '
$ws.Cells.Item(36, 4).Value = '3.712.50'
$ws.Cells.Item(36, 5).Value = '  -0.01%  '
$ws.Cells.Item(37, 4).Value = '3.35'
$ws.Cells.Item(37, 5).Value = '  -0.51%  '
$ws.Cells.Item(38, 4).Value = '0.137'
$ws.Cells.Item(38, 5).Value = '  +4.89%  '
$ws.Cells.Item(39, 4).Value = '34.85'
$ws.Cells.Item(39, 5).Value = '  +1.14%  '
$ws.Cells.Item(40, 5).Value = '  +1.71%  '
$ws.Cells.Item(41, 5).Value = '  -0.81%  '
$ws.Cells.Item(42, 4).Value = '0.0₃0674'
$ws.Cells.Item(42, 5).Value = '  +0.39%  '
$ws.Cells.Item(43, 4).Value = '0.335'
$ws.Cells.Item(43, 5).Value = '  +0.10%  '
$ws.Cells.Item(44, 4).Value = '3.25'
$ws.Cells.Item(44, 5).Value = '  -1.37%  '
$ws.Cells.Item(45, 5).Value = '  +1.15%  '
$ws.Cells.Item(47, 5).Value = '  +0.30%  '
$ws.Cells.Item(48, 5).Value = '  -0.14%  '
$ws.Cells.Item(49, 5).Value = '  -1.00%  '
$ws.Cells.Item(50, 5).Value = '  +2.80%  '
$ws.Cells.Item(51, 5).Value = '  -1.69%  '

# Restore the default (unstyled) cell style now that the text value is set,
# so no extra style gets attached to these cells.
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).Style = "Normal"
}
